$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.739.01"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.598.94"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.65"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.73"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.824.02"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.582.15"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.04"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "26.710.51"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "0.0₃0743"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.69"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.16"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.15"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0508"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  +17.28%  "
$ws.Range("D35").Value = "1.276.88"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.595"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.44"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.63"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "1.735.22"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.41"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("E51").Value = "  +0.26%  "
